$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 4226.5674
$ws.Range("I32").Value = 4180.143
$ws.Range("J32").Value = 4237.4
$ws.Range("K32").Value = 4180.143
$ws.Range("L32").Value = 4237.4
$ws.Range("M32").Value = -3854.143
$ws.Range("N32").Value = -4889.4
$ws.Range("H51").Value = 7000
$ws.Range("I51").Value = 7000
$ws.Range("K51").Value = 7000
$ws.Range("M51").Value = -6516
$ws.Range("H58").Value = 4637
$ws.Range("I58").Value = 2164.111
$ws.Range("K58").Value = 6492.333
$ws.Range("M58").Value = -6342.333
$ws.Range("H64").Value = 12957.692
$ws.Range("I64").Value = 20666.666
$ws.Range("J64").Value = 6350
$ws.Range("K64").Value = 20666.666
$ws.Range("L64").Value = 6350
$ws.Range("M64").Value = -20418.666
$ws.Range("N64").Value = -6846
$ws.Range("H67").Value = 12957.692
$ws.Range("I67").Value = 20666.666
$ws.Range("J67").Value = 6350
$ws.Range("K67").Value = 20666.666
$ws.Range("L67").Value = 6350
$ws.Range("M67").Value = -19808.666
$ws.Range("N67").Value = -8066
$ws.Range("H132").Value = 71437870
$ws.Range("I132").Value = 142870910
$ws.Range("K132").Value = 428612730
$ws.Range("M132").Value = -428610200
$ws.Range("H138").Value = 5251.8335
$ws.Range("J138").Value = 5155.222
$ws.Range("L138").Value = 15465.666
$ws.Range("N138").Value = -25745.666
$ws.Range("H141").Value = 4676.8
$ws.Range("I141").Value = 4940.857
$ws.Range("K141").Value = 14822.571
$ws.Range("M141").Value = -9642.571

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4774.0195
$ws.Range("I32").Value = 3301.976
$ws.Range("J32").Value = 11643.556
$ws.Range("K32").Value = 3301.976
$ws.Range("L32").Value = 11643.556
$ws.Range("M32").Value = -3014.976
$ws.Range("N32").Value = -12217.556
$ws.Range("H45").Value = 11993418
$ws.Range("I45").Value = 23977588
$ws.Range("J45").Value = 9248.833000000001
$ws.Range("K45").Value = 23977588
$ws.Range("L45").Value = 9248.833000000001
$ws.Range("M45").Value = -23977211
$ws.Range("N45").Value = -10002.833
$ws.Range("H61").Value = 7318.75
$ws.Range("I61").Value = 12688.333
$ws.Range("K61").Value = 12688.333
$ws.Range("M61").Value = -12476.333
$ws.Range("H132").Value = 18471.375
$ws.Range("I132").Value = 23156.2
$ws.Range("K132").Value = 69468.60000000001
$ws.Range("M132").Value = -66938.60000000001
$ws.Range("H136").Value = 7318.75
$ws.Range("I136").Value = 12688.333
$ws.Range("K136").Value = 38064.999
$ws.Range("M136").Value = -35514.999
$ws.Range("H137").Value = 50780
$ws.Range("J137").Value = 50780
$ws.Range("L137").Value = 50780
$ws.Range("N137").Value = -60980

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H54").Value = 8360.286
$ws.Range("I54").Value = 9704.799999999999
$ws.Range("J54").Value = 4999
$ws.Range("K54").Value = 9704.799999999999
$ws.Range("L54").Value = 4999
$ws.Range("M54").Value = -9220.799999999999
$ws.Range("N54").Value = -5967
$ws.Range("H58").Value = 49389.5
$ws.Range("J58").Value = 49389.5
$ws.Range("L58").Value = 49389.5
$ws.Range("N58").Value = -49977.5
$ws.Range("H106").Value = 25595.75
$ws.Range("J106").Value = 25595.75
$ws.Range("L106").Value = 25595.75
$ws.Range("N106").Value = -28119.75
$ws.Range("H134").Value = 11367.407
$ws.Range("I134").Value = 11924.45
$ws.Range("K134").Value = 35773.35000000001
$ws.Range("M134").Value = -33238.35000000001

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H53").Value = 39665
$ws.Range("J53").Value = 39665
$ws.Range("L53").Value = 39665
$ws.Range("N53").Value = -40879
$ws.Range("H58").Value = 2428.6316
$ws.Range("I58").Value = 2236.5454
$ws.Range("K58").Value = 2236.5454
$ws.Range("M58").Value = -2033.5454
$ws.Range("H99").Value = 4614.154
$ws.Range("I99").Value = 4416.6665
$ws.Range("K99").Value = 4416.6665
$ws.Range("M99").Value = -2918.6665
$ws.Range("H126").Value = 4614.154
$ws.Range("I126").Value = 4416.6665
$ws.Range("K126").Value = 13249.9995
$ws.Range("M126").Value = -10779.9995
$ws.Range("H132").Value = 115996.445
$ws.Range("I132").Value = 169162.33
$ws.Range("J132").Value = 9664.666999999999
$ws.Range("K132").Value = 507486.99
$ws.Range("L132").Value = 28994.001
$ws.Range("M132").Value = -504956.99
$ws.Range("N132").Value = -34054.001
$ws.Range("H136").Value = 2428.6316
$ws.Range("I136").Value = 2236.5454
$ws.Range("K136").Value = 6709.6362
$ws.Range("M136").Value = -4159.6362

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 1653.6364
$ws.Range("I3").Value = 1148.75
$ws.Range("J3").Value = 3000
$ws.Range("K3").Value = 3446.25
$ws.Range("L3").Value = 9000
$ws.Range("M3").Value = -3334.25
$ws.Range("N3").Value = -9224
$ws.Range("H64").Value = 0
$ws.Range("I64").Value = 0
$ws.Range("K64").Value = 0
$ws.Range("M64").Value = ""
$ws.Range("H67").Value = 0
$ws.Range("I67").Value = 0
$ws.Range("K67").Value = 0
$ws.Range("M67").Value = ""
$ws.Range("H107").Value = 784.1667
$ws.Range("J107").Value = 967.9167
$ws.Range("L107").Value = 2903.7501
$ws.Range("N107").Value = -6743.7501
$ws.Range("H113").Value = 2275.4
$ws.Range("J113").Value = 1779.44
$ws.Range("L113").Value = 5338.32
$ws.Range("N113").Value = -9678.32
$ws.Range("H131").Value = 3045.9546
$ws.Range("J131").Value = 5909.375
$ws.Range("L131").Value = 17728.125
$ws.Range("N131").Value = -27808.125
$ws.Range("H137").Value = 3262.6296
$ws.Range("I137").Value = 2257.1667
$ws.Range("J137").Value = 4067
$ws.Range("K137").Value = 6771.500100000001
$ws.Range("L137").Value = 12201
$ws.Range("M137").Value = -1671.500100000001
$ws.Range("N137").Value = -22401

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 37.875
$ws.Range("J2").Value = 9.142858
$ws.Range("L2").Value = 9.142858
$ws.Range("N2").Value = -235.142858
$ws.Range("H70").Value = 18187008
$ws.Range("I70").Value = 22227120
$ws.Range("J70").Value = 6500
$ws.Range("K70").Value = 22227120
$ws.Range("L70").Value = 6500
$ws.Range("M70").Value = -22226850
$ws.Range("N70").Value = -7040
$ws.Range("H73").Value = 18187008
$ws.Range("I73").Value = 22227120
$ws.Range("J73").Value = 6500
$ws.Range("K73").Value = 22227120
$ws.Range("L73").Value = 6500
$ws.Range("M73").Value = -22226184
$ws.Range("N73").Value = -8372
$ws.Range("H110").Value = 75351
$ws.Range("J110").Value = 75351
$ws.Range("L110").Value = 75351
$ws.Range("N110").Value = -83531
$ws.Range("H135").Value = 76000
$ws.Range("J135").Value = 76000
$ws.Range("L135").Value = 76000
$ws.Range("N135").Value = -86140

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 44279.57
$ws.Range("J22").Value = 2270.2144
$ws.Range("L22").Value = 2270.2144
$ws.Range("N22").Value = -2860.2144
$ws.Range("H27").Value = 44279.57
$ws.Range("J27").Value = 2270.2144
$ws.Range("L27").Value = 2270.2144
$ws.Range("N27").Value = -2484.2144
$ws.Range("H68").Value = 2477.25
$ws.Range("I68").Value = 2363.8
$ws.Range("K68").Value = 2363.8
$ws.Range("M68").Value = -1614.8
$ws.Range("H71").Value = 2477.25
$ws.Range("I71").Value = 2363.8
$ws.Range("K71").Value = 11819
$ws.Range("M71").Value = -8075
$ws.Range("H100").Value = 4328.2856
$ws.Range("I100").Value = 3968.4375
$ws.Range("J100").Value = 5479.8
$ws.Range("K100").Value = 3968.4375
$ws.Range("L100").Value = 5479.8
$ws.Range("M100").Value = -3427.4375
$ws.Range("N100").Value = -6561.8
$ws.Range("H136").Value = 54902.45
$ws.Range("J136").Value = 7999.5
$ws.Range("L136").Value = 23998.5
$ws.Range("N136").Value = -29098.5

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H28").Value = 19000.5
$ws.Range("J28").Value = 19000.5
$ws.Range("L28").Value = 19000.5
$ws.Range("N28").Value = -19696.5
$ws.Range("H29").Value = 10835.167
$ws.Range("I29").Value = 10000
$ws.Range("K29").Value = 10000
$ws.Range("M29").Value = -9710
$ws.Range("H31").Value = 11836.333
$ws.Range("J31").Value = 9603.799999999999
$ws.Range("L31").Value = 9603.799999999999
$ws.Range("N31").Value = -10299.8
$ws.Range("H39").Value = 28500
$ws.Range("I39").Value = 28500
$ws.Range("K39").Value = 28500
$ws.Range("M39").Value = -28087
$ws.Range("H132").Value = 21983304
$ws.Range("I132").Value = 29418078
$ws.Range("J132").Value = 918113.5600000001
$ws.Range("K132").Value = 88254234
$ws.Range("L132").Value = 2754340.68
$ws.Range("M132").Value = -88251704
$ws.Range("N132").Value = -2759400.68
$ws.Range("H136").Value = 4830.0786
$ws.Range("I136").Value = 5786.1177
$ws.Range("J136").Value = 2918
$ws.Range("K136").Value = 17358.3531
$ws.Range("L136").Value = 8754
$ws.Range("M136").Value = -14808.3531
$ws.Range("N136").Value = -13854
